$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.019.60'
$ws.Range("E2").Value = '  -1.10%  '
$ws.Range("D3").Value = '1.822.53'
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.52'
$ws.Range("E5").Value = '  -1.69%  '
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4638'
$ws.Range("E7").Value = '  -2.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3660'
$ws.Range("E8").Value = '  -0.82%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07236'
$ws.Range("E9").Value = '  -2.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8590'
$ws.Range("E11").Value = '  -3.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07642'
$ws.Range("E12").Value = '  +4.18%  '
$ws.Range("D13").Value = '1.780.64'
$ws.Range("E13").Value = '  -5.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.327'
$ws.Range("E14").Value = '  -2.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.69'
$ws.Range("E16").Value = '  -2.09%  '
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("E18").Value = '  -1.77%  '
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").Value = '27.201.66'
$ws.Range("E20").Value = '  -1.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.47'
$ws.Range("E21").Value = '  -2.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.148'
$ws.Range("E22").Value = '  -2.70%  '
$ws.Range("E23").Value = '  -1.18%  '
$ws.Range("D24").Value = '1.979.31'
$ws.Range("E24").Value = '  -5.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.78'
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.839'
$ws.Range("E26").Value = '  -2.96%  '
$ws.Range("E27").Value = '  -2.69%  '
$ws.Range("E28").Value = '  -4.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.100'
$ws.Range("E29").Value = '  -2.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.30'
$ws.Range("E30").Value = '  -1.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08836'
$ws.Range("E31").Value = '  -1.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.427'
$ws.Range("E34").Value = '  -4.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7194'
$ws.Range("E35").Value = '  -4.09%  '
$ws.Range("E36").Value = '  -2.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05250'
$ws.Range("E37").Value = '  -1.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01926'
$ws.Range("E38").Value = '  -1.48%  '
$ws.Range("E39").Value = '  +1.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.929'
$ws.Range("E40").Value = '  -1.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.134'
$ws.Range("E41").Value = '  -1.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5170'
$ws.Range("E42").Value = '  -2.69%  '
$ws.Range("E43").Value = '  -2.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8589'
$ws.Range("E44").Value = '  -15.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.156'
$ws.Range("E46").Value = '  -2.47%  '
$ws.Range("E47").Value = '  -0.34%  '
$ws.Range("E48").Value = '  -3.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '102.69'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06254'
$ws.Range("E50").Value = '  -0.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.617'
$ws.Range("E51").Value = '  -3.32%  '
